$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new application row (row 6).
# The tech-stack string (B6) is written first so it lands earlier in the
# workbook's shared-strings table, matching the target ordering.
$ws.Range("B6").Value = "C#, .NET Core, MVC, Web API, REST, ADO.NET, Nunit testing, Moq, Nlog, Ocelot API Gateway, Exception handeling, Microservices"
$ws.Range("A6").Value = "MicroservicesExperiments\Geodesics"

# Widen column A to fit the new, longer content
$ws.Columns.Item(1).ColumnWidth = 34.333333333333336

# Update selection to match the target workbook state
$ws.Range("B8").Select()
